$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 1.02
    "C2" = 1.038637502042585
    "D2" = 1.048689865595021
    "E2" = 1.047114037982982
    "F2" = 1.057644199255419
    "I2" = 1.033448618331143
    "J2" = 1.043733857995824
    "K2" = 1.051449008593891
    "L2" = 1.049877581866902
    "M2" = 1.060378604807207
    "N2" = 1.018486864579258
    "B3" = 1.02
    "C3" = 1.039557447109204
    "D3" = 1.049553471098263
    "E3" = 1.047943685838606
    "F3" = 1.058583732764481
    "I3" = 1.033536788815871
    "J3" = 1.04429882646677
    "K3" = 1.052124576786032
    "L3" = 1.050518962222228
    "M3" = 1.061131696044088
    "N3" = 1.018676380304725
    "B4" = 1.02
    "C4" = 1.040153242855415
    "D4" = 1.050113136007309
    "E4" = 1.048481406810249
    "F4" = 1.059192773332784
    "I4" = 1.033592421776112
    "J4" = 1.044664294445915
    "K4" = 1.052561926259448
    "L4" = 1.050934211755989
    "M4" = 1.061619466944717
    "N4" = 1.018798917704992
    "B5" = 1.02
    "C5" = 1.040403840949444
    "D5" = 1.05034862210182
    "E5" = 1.048707674610438
    "F5" = 1.059449075611044
    "I5" = 1.033615469519347
    "J5" = 1.044817911038616
    "K5" = 1.052745837342461
    "L5" = 1.051108837312064
    "M5" = 1.061824637070426
    "N5" = 1.018850409976366
    "B6" = 1.02
    "C6" = 1.040445924798418
    "D6" = 1.050388173109466
    "E6" = 1.048745678232994
    "F6" = 1.059492125186162
    "I6" = 1.033619319361386
    "J6" = 1.044843702373793
    "K6" = 1.052776719679168
    "L6" = 1.051138160879472
    "M6" = 1.06185909251968
    "N6" = 1.018859054425947
    "B7" = 1.02
    "C7" = 1.040156590867619
    "D7" = 1.050116281788491
    "E7" = 1.048484429388861
    "F7" = 1.059196197030395
    "I7" = 1.033592731079271
    "J7" = 1.044666347182527
    "K7" = 1.052564383496019
    "L7" = 1.050936544898055
    "M7" = 1.061622208002646
    "N7" = 1.018799605836014
    "B8" = 1.02
    "C8" = 1.038948291329484
    "D8" = 1.048981547259697
    "E8" = 1.047394237517841
    "F8" = 1.057961490360973
    "I8" = 1.033478709470654
    "J8" = 1.043924812659599
    "K8" = 1.051677275408874
    "L8" = 1.050094290242942
    "M8" = 1.060633017250438
    "N8" = 1.018550931080006
    "B9" = 1.02
    "C9" = 1.036823226591152
    "D9" = 1.046988613476559
    "E9" = 1.045480014686694
    "F9" = 1.055794272146152
    "I9" = 1.03326694496509
    "J9" = 1.042617385401031
    "K9" = 1.050115756790031
    "L9" = 1.048611974426488
    "M9" = 1.058893602299941
    "N9" = 1.018112049585399
    "B10" = 1.02
    "C10" = 1.03540935702664
    "D10" = 1.045664531035805
    "E10" = 1.044208554287813
    "F10" = 1.054355265706302
    "I10" = 1.033118512475238
    "J10" = 1.041745332053282
    "K10" = 1.04907595106981
    "L10" = 1.047625078768965
    "M10" = 1.057736539940333
    "N10" = 1.017819027521553
    "B11" = 1.02
    "C11" = 1.034797825016728
    "D11" = 1.045092285161297
    "E11" = 1.043659129361148
    "F11" = 1.053733557093965
    "I11" = 1.033052525745116
    "J11" = 1.04136763435228
    "K11" = 1.048626006170927
    "L11" = 1.047198069563435
    "M11" = 1.057236140804295
    "N11" = 1.017692048142487
    "B12" = 1.02
    "C12" = 1.034570778597212
    "D12" = 1.044879892975176
    "E12" = 1.043455219183885
    "F12" = 1.053502837368273
    "I12" = 1.033027758312509
    "J12" = 1.041227327786687
    "K12" = 1.048458922639204
    "L12" = 1.047039509121646
    "M12" = 1.057050364291163
    "N12" = 1.017644867921616
    "B13" = 1.02
    "C13" = 1.03461947610843
    "D13" = 1.044925444309731
    "E13" = 1.043498950859312
    "F13" = 1.053552317979771
    "I13" = 1.033033082632737
    "J13" = 1.041257424583433
    "K13" = 1.048494760519915
    "L13" = 1.047073518601594
    "M13" = 1.057090209707652
    "N13" = 1.017654988885017
    "B14" = 1.02
    "C14" = 1.034779055164424
    "D14" = 1.045074725373542
    "E14" = 1.043642270590947
    "F14" = 1.053714481425359
    "I14" = 1.033050483701368
    "J14" = 1.041356036822855
    "K14" = 1.048612194047456
    "L14" = 1.047184961884964
    "M14" = 1.05722078253406
    "N14" = 1.017688148502496
    "B15" = 1.02
    "C15" = 1.034877390877982
    "D15" = 1.045166724329275
    "E15" = 1.043730597274382
    "F15" = 1.053814423622053
    "I15" = 1.03306117103274
    "J15" = 1.041416793412794
    "K15" = 1.048684554858602
    "L15" = 1.047253632404773
    "M15" = 1.057301245260014
    "N15" = 1.01770857733859
    "B16" = 1.02
    "C16" = 1.035449956881823
    "D16" = 1.045702532220409
    "E16" = 1.044245041663738
    "F16" = 1.054396555895915
    "I16" = 1.033122855722693
    "J16" = 1.041770396745766
    "K16" = 1.049105818816168
    "L16" = 1.047653424877664
    "M16" = 1.057769762890049
    "N16" = 1.017827452678059
    "B17" = 1.02
    "C17" = 1.035809296076712
    "D17" = 1.04603892343433
    "E17" = 1.044568041553799
    "F17" = 1.054762085544376
    "I17" = 1.033161090277597
    "J17" = 1.041992178711957
    "K17" = 1.049370147218027
    "L17" = 1.047904291521495
    "M17" = 1.058063817652875
    "N17" = 1.017901993866519
    "B18" = 1.02
    "C18" = 1.036018958352609
    "D18" = 1.046235239856255
    "E18" = 1.044756550445125
    "F18" = 1.054975426876263
    "I18" = 1.033183226386407
    "J18" = 1.042121531346614
    "K18" = 1.04952435415337
    "L18" = 1.048050648884126
    "M18" = 1.058235394211575
    "N18" = 1.017945462910054
    "B19" = 1.02
    "C19" = 1.036090458895177
    "D19" = 1.046302196465064
    "E19" = 1.044820845464928
    "F19" = 1.055048193449461
    "I19" = 1.033190746145555
    "J19" = 1.042165635685675
    "K19" = 1.04957693955365
    "L19" = 1.048100558198759
    "M19" = 1.058293907415733
    "N19" = 1.017960283094091
    "B20" = 1.02
    "C20" = 1.035770735566785
    "D20" = 1.046002820956465
    "E20" = 1.044533375486392
    "F20" = 1.054722853793137
    "I20" = 1.03315700518199
    "J20" = 1.041968384532184
    "K20" = 1.049341784309287
    "L20" = 1.047877372679641
    "M20" = 1.05803226220668
    "N20" = 1.017893997289255
    "B21" = 1.02
    "C21" = 1.034732060234796
    "D21" = 1.045030761242248
    "E21" = 1.043600061777328
    "F21" = 1.053666722507799
    "I21" = 1.033045366610971
    "J21" = 1.041326998321101
    "K21" = 1.048577611519321
    "L21" = 1.047152143231232
    "M21" = 1.057182329484974
    "N21" = 1.017678384216386
    "B22" = 1.02
    "C22" = 1.034079604969926
    "D22" = 1.044420546511431
    "E22" = 1.04301423881739
    "F22" = 1.053003909732475
    "I22" = 1.032973688231601
    "J22" = 1.040923659725181
    "K22" = 1.048097412754599
    "L22" = 1.046696451376756
    "M22" = 1.056648487626106
    "N22" = 1.017542736470518
    "B23" = 1.02
    "C23" = 1.034425426464853
    "D23" = 1.044743941574979
    "E23" = 1.043324700440835
    "F23" = 1.053355163254467
    "I23" = 1.033011827039513
    "J23" = 1.041137483749014
    "K23" = 1.048351949465987
    "L23" = 1.046937994471496
    "M23" = 1.056931435214154
    "N23" = 1.017614653654141
    "B24" = 1.02
    "C24" = 1.035788159198185
    "D24" = 1.046019133787571
    "E24" = 1.044549039254933
    "F24" = 1.054740580518722
    "I24" = 1.033158851572156
    "J24" = 1.041979136125842
    "K24" = 1.049354600198487
    "L24" = 1.047889536049622
    "M24" = 1.058046520571164
    "N24" = 1.01789761062793
    "B25" = 1.02
    "C25" = 1.037372111626546
    "D25" = 1.047503041262554
    "E25" = 1.045974067661563
    "F25" = 1.05635353371841
    "I25" = 1.033322972650403
    "J25" = 1.042955468384112
    "K25" = 1.050519239786542
    "L25" = 1.048994962312543
    "M25" = 1.059342839812476
    "N25" = 1.018225589442042
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}